$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 304, shifting existing rows 304-314 down to 305-315.
$ws.Rows.Item(304).Insert()

# Populate the newly inserted row 304 with the new weekly data point.
$ws.Cells.Item(304, 1).Value = 4
$ws.Cells.Item(304, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(304, 3).Value = "Los Lagos"
$ws.Cells.Item(304, 4).Value = 44753
$ws.Cells.Item(304, 5).Value = 10
$ws.Cells.Item(304, 6).Value = 100112040
$ws.Cells.Item(304, 7).Value = "Cilantro"
$ws.Cells.Item(304, 8).Value = "Sin especificar"
$ws.Cells.Item(304, 9).Value = "Primera"
$ws.Cells.Item(304, 10).Value = 70
$ws.Cells.Item(304, 11).Value = 14000
$ws.Cells.Item(304, 12).Value = 14000
$ws.Cells.Item(304, 13).Value = 14000
$ws.Cells.Item(304, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(304, 15).Value = "Región Metropolitana"
$ws.Cells.Item(304, 16).Value = 389
$ws.Cells.Item(304, 17).Value = 36
$ws.Cells.Item(304, 18).Value = "Hortaliza"
